$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Rose"
$ws.Range("B3").Value = "Mathew"
$ws.Range("A4").Value = "Vishnupriyan"
$ws.Range("B4").Value = "V S"
